$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The numeric-looking IDs ("32", "1", "4", "9009", "26") must be stored as
# text, not numbers, so force text format before typing them in (otherwise
# Excel auto-coerces digit strings to numeric cells).
$dataRange = $ws.Range("A2:H2")
$dataRange.NumberFormat = "@"

$ws.Range("A2").Value = "GenCor2019"
$ws.Range("B2").Value = "AR-X"
$ws.Range("C2").Value = "32"
$ws.Range("D2").Value = "1"
$ws.Range("E2").Value = "4"
$ws.Range("F2").Value = "1"
$ws.Range("G2").Value = "9009"
$ws.Range("H2").Value = "26"

# Restore the default (unformatted) style on these cells so they don't carry
# an explicit text-number-format style index.
$dataRange.Style = "Normal"

# Standard Excel page margins (inches -> points).
$ws.PageSetup.LeftMargin = 0.7 * 72
$ws.PageSetup.RightMargin = 0.7 * 72
$ws.PageSetup.TopMargin = 0.75 * 72
$ws.PageSetup.BottomMargin = 0.75 * 72
$ws.PageSetup.HeaderMargin = 0.3 * 72
$ws.PageSetup.FooterMargin = 0.3 * 72
